$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.00"
$ws.Range("D3").Value = "'22.11"
$ws.Range("D4").Value = "'5.302"
$ws.Range("D5").Value = "'0.05874"
$ws.Range("D7").Value = "'6.379"
$ws.Range("D8").Value = "'0.8158"
$ws.Range("D9").Value = "'0.9608"
$ws.Range("D10").Value = "'0.1417"
$ws.Range("D11").Value = "'0.03550"
$ws.Range("D12").Value = "'0.07327"
$ws.Range("D13").Value = "'0.03041"
$ws.Range("D14").Value = "'4.429"
$ws.Range("D15").Value = "'0.09388"
$ws.Range("D16").Value = "'0.001603"
$ws.Range("D17").Value = "'0.04821"
$ws.Range("D18").Value = "'0.0005903"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("D19").Value = "'0.006000"
$ws.Range("D20").Value = "'0.004080"
$ws.Range("D21").Value = "'0.0009866"
$ws.Range("D22").Value = "'0.00009705"
$ws.Range("D24").Value = "'2.180"
$ws.Range("D25").Value = "'0.3259"
$ws.Range("D40").Value = "'0.03859"
$ws.Range("D41").Value = "'0.006613"
$ws.Range("D42").Value = "'0.1073"
$ws.Range("D43").Value = "'0.003001"
$ws.Range("D44").Value = "'0.005298"
$ws.Range("D45").Value = "'0.00005666"
$ws.Range("D47").Value = "'0.7753"
$ws.Range("D48").Value = "'0.08639"
$ws.Range("E48").Value = "47BOLOBOLO"
$ws.Range("D50").Value = "'0.01010"
